$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-07-22 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-23 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("114÷7=16, 2", $true, $false, $false, $false, $false, $true, 1, $false, "853÷7=121, 6", 2) | Out-Null
$d.Content.Find.Execute("967÷5=193, 2", $true, $false, $false, $false, $false, $true, 1, $false, "778÷8=97, 2", 2) | Out-Null
$d.Content.Find.Execute("331÷5=66, 1", $true, $false, $false, $false, $false, $true, 1, $false, "187÷4=46, 3", 2) | Out-Null
$d.Content.Find.Execute("554÷9=61, 5", $true, $false, $false, $false, $false, $true, 1, $false, "282÷6=47, 0", 2) | Out-Null
$d.Content.Find.Execute("915÷6=152, 3", $true, $false, $false, $false, $false, $true, 1, $false, "959÷8=119, 7", 2) | Out-Null
$d.Content.Find.Execute("362÷5=72, 2", $true, $false, $false, $false, $false, $true, 1, $false, "822÷4=205, 2", 2) | Out-Null
$d.Content.Find.Execute("808÷8=101, 0", $true, $false, $false, $false, $false, $true, 1, $false, "488÷2=244, 0", 2) | Out-Null
$d.Content.Find.Execute("955÷8=119, 3", $true, $false, $false, $false, $false, $true, 1, $false, "603÷3=201, 0", 2) | Out-Null
$d.Content.Find.Execute("576÷5=115, 1", $true, $false, $false, $false, $false, $true, 1, $false, "720÷4=180, 0", 2) | Out-Null
$d.Content.Find.Execute("920÷7=131, 3", $true, $false, $false, $false, $false, $true, 1, $false, "627÷3=209, 0", 2) | Out-Null
$d.Content.Find.Execute("201÷8=25, 1", $true, $false, $false, $false, $false, $true, 1, $false, "689÷4=172, 1", 2) | Out-Null
$d.Content.Find.Execute("420÷8=52, 4", $true, $false, $false, $false, $false, $true, 1, $false, "445÷8=55, 5", 2) | Out-Null
$d.Content.Find.Execute("455÷3=151, 2", $true, $false, $false, $false, $false, $true, 1, $false, "779÷5=155, 4", 2) | Out-Null
$d.Content.Find.Execute("123÷8=15, 3", $true, $false, $false, $false, $false, $true, 1, $false, "274÷2=137, 0", 2) | Out-Null
$d.Content.Find.Execute("650÷2=325, 0", $true, $false, $false, $false, $false, $true, 1, $false, "798÷3=266, 0", 2) | Out-Null
$d.Content.Find.Execute("330÷2=165, 0", $true, $false, $false, $false, $false, $true, 1, $false, "300÷5=60, 0", 2) | Out-Null
$d.Content.Find.Execute("182÷2=91, 0", $true, $false, $false, $false, $false, $true, 1, $false, "437÷6=72, 5", 2) | Out-Null
$d.Content.Find.Execute("782÷6=130, 2", $true, $false, $false, $false, $false, $true, 1, $false, "231÷8=28, 7", 2) | Out-Null
$d.Content.Find.Execute("695÷7=99, 2", $true, $false, $false, $false, $false, $true, 1, $false, "661÷6=110, 1", 2) | Out-Null
$d.Content.Find.Execute("125÷4=31, 1", $true, $false, $false, $false, $false, $true, 1, $false, "687÷6=114, 3", 2) | Out-Null
$d.Content.Find.Execute("739÷8=92, 3", $true, $false, $false, $false, $false, $true, 1, $false, "922÷7=131, 5", 2) | Out-Null
$d.Content.Find.Execute("948÷3=316, 0", $true, $false, $false, $false, $false, $true, 1, $false, "608÷2=304, 0", 2) | Out-Null
$d.Content.Find.Execute("571÷5=114, 1", $true, $false, $false, $false, $false, $true, 1, $false, "343÷3=114, 1", 2) | Out-Null
$d.Content.Find.Execute("377÷9=41, 8", $true, $false, $false, $false, $false, $true, 1, $false, "957÷5=191, 2", 2) | Out-Null
$d.Content.Find.Execute("315÷8=39, 3", $true, $false, $false, $false, $false, $true, 1, $false, "954÷6=159, 0", 2) | Out-Null
